# Generate Report for Handback
#
# Replaces the two file "tokens" used throughout the handback-status report:
#   1c0696b1-d9a5-4d71-b2b8-029e60c0a26b  ->  5a04021b-a533-40d8-9da6-7aaf68baea15
#   3019a9a2-7081-43f1-ba05-ee32a65e7bf0  ->  fffffcdd500d-9591-4e41-ae7a-bcfc200ac221
# along with their hashed xliff file names and the associated timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.md"
$wsOverview.Range("B2").Value = "e2e\5a04021b-a533-40d8-9da6-7aaf68baea15.md"
$wsOverview.Range("G2").Value = "2016-08-25 00:59:06"

$wsOverview.Range("A3").Value = "fffffcdd500d-9591-4e41-ae7a-bcfc200ac221.md"
$wsOverview.Range("B3").Value = "e2e\fffffcdd500d-9591-4e41-ae7a-bcfc200ac221.md"
$wsOverview.Range("G3").Value = "2016-08-25 00:59:06"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("A2").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.md"
$wsZhCn.Range("G2").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-25 00:58:57"
$wsZhCn.Range("I2").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.md"
$wsZhCn.Range("J2").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-25 00:59:27"

$wsZhCn.Range("A3").Value = "fffffcdd500d-9591-4e41-ae7a-bcfc200ac221.md"
$wsZhCn.Range("G3").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-25 00:58:57"
$wsZhCn.Range("I3").Value = "fffffcdd500d-9591-4e41-ae7a-bcfc200ac221.md"
$wsZhCn.Range("J3").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-25 00:59:27"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe.Range("A2").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.md"
$wsDeDe.Range("G2").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-25 00:59:06"
$wsDeDe.Range("I2").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.md"
$wsDeDe.Range("J2").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-25 00:59:34"

$wsDeDe.Range("A3").Value = "fffffcdd500d-9591-4e41-ae7a-bcfc200ac221.md"
$wsDeDe.Range("G3").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-25 00:59:06"
$wsDeDe.Range("I3").Value = "fffffcdd500d-9591-4e41-ae7a-bcfc200ac221.md"
$wsDeDe.Range("J3").Value = "5a04021b-a533-40d8-9da6-7aaf68baea15.70ceb74897fb9812e6766e5e58af29a09d2a00f7.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-25 00:59:34"

# ---------------------------------------------------------------------------
# Update the hyperlink display text (the targets / relationship ids are left
# untouched) to reflect the renamed files. Capturing each Hyperlink object
# via foreach and mutating TextToDisplay updates it in-place instead of
# creating a brand new hyperlink entry.
# ---------------------------------------------------------------------------
function Update-HyperlinkDisplay($worksheet, $replacements) {
    $links = @()
    foreach ($hl in $worksheet.Hyperlinks) {
        $links += $hl
    }
    foreach ($hl in $links) {
        $text = $hl.TextToDisplay
        foreach ($pair in $replacements) {
            $text = $text -replace [Regex]::Escape($pair[0]), $pair[1]
        }
        $hl.TextToDisplay = $text
    }
}

$replacements = @(
    , @("1c0696b1-d9a5-4d71-b2b8-029e60c0a26b", "5a04021b-a533-40d8-9da6-7aaf68baea15")
    , @("3019a9a2-7081-43f1-ba05-ee32a65e7bf0", "fffffcdd500d-9591-4e41-ae7a-bcfc200ac221")
)

Update-HyperlinkDisplay $wsOverview $replacements
Update-HyperlinkDisplay $wsZhCn $replacements
Update-HyperlinkDisplay $wsDeDe $replacements
